$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "last updated" timestamp string (cell A1)
$ws.Range("A1").Value = "Datos actualizados a 29 de Abril de 2020 a las 10:22"

# Row 27 - Israel
$ws.Range("B27").Value = 15782
$ws.Range("C27").Value = 54
$ws.Range("D27").Value = 7929
$ws.Range("E27").Value = 7641
$ws.Range("F27").Value = 120
$ws.Range("G27").Value = 2
$ws.Range("H27").Value = 212

# Row 33 - Polonia
$ws.Range("B33").Value = 12415
$ws.Range("C33").Value = 197
$ws.Range("E33").Value = 8784
$ws.Range("G33").Value = 10
$ws.Range("H33").Value = 606

# Row 43 - Filipinas
$ws.Range("B43").Value = 8212
$ws.Range("C43").Value = 254
$ws.Range("D43").Value = 1023
$ws.Range("E43").Value = 6631
$ws.Range("G43").Value = 28
$ws.Range("H43").Value = 558

# Row 84 - Eslovaquia
$ws.Range("B84").Value = 1391
$ws.Range("C84").Value = 7
$ws.Range("D84").Value = 484
$ws.Range("E84").Value = 885
$ws.Range("F84").Value = 8
$ws.Range("G84").Value = 2
$ws.Range("H84").Value = 22

# Row 136 - Birmania
$ws.Range("E136").Value = 128
$ws.Range("G136").Value = 1
$ws.Range("H136").Value = 6
